# deposits view shows forecast
#
# The deposit-estimation "forecast" block (F36:F41) on the active sheet
# gets relabeled (now typed with explicit "Period"/"decimal" prefixes,
# matching the updated DepositEstimations data-contract), the helper
# DepositRateLine property list that used to live in column K is moved to
# column H, the F column is widened to fit the new (longer) labels, the
# big "DepositCalculationData" class documentation text box is removed
# from the drawing layer, and the saved cursor position moves to I25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relabel the forecast rows (F36:F41) -------------------------------
# Old:   ProcentsInThisMonth / CurrencyRateOnThisMonthPayment / ProcentsUpToFinish /
#        CurrencyRateOnFinish / DevaluationInUsd / ProfitInUsd
# New:   Period PeriodForThisMonthPayment / decimal ProcentsInThisMonth /
#        Period PeriodForUpToEndPayment / decimal ProcentsUpToFinish /
#        decimal DevaluationInUsd / decimal ProfitInUsd
$ws.Range("F36").Value = "Period PeriodForThisMonthPayment"
$ws.Range("F37").Value = "decimal ProcentsInThisMonth"
$ws.Range("F38").Value = "Period PeriodForUpToEndPayment"
$ws.Range("F39").Value = "decimal ProcentsUpToFinish"
$ws.Range("F40").Value = "decimal DevaluationInUsd"
$ws.Range("F41").Value = "decimal ProfitInUsd"

# --- 2. Move the DepositRateLine property list from column K to column H -
$ws.Range("K15:K20").Copy($ws.Range("H15:H20"))
$ws.Range("K15:K20").Clear()

# --- 3. Widen column F to fit the longer labels ---------------------------
$ws.Columns("F").ColumnWidth = 32.8

# --- 4. Remove the big "DepositCalculationData" comment text box ----------
$ws.Shapes.Item("TextBox 1").Delete()

# --- 5. Move the saved selection to I25 ------------------------------------
$ws.Range("I25").Select()
